# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" column (E16:E55) used to list the 40 monthly periods
# in descending order (2003 down to 1612). This update re-lists them in
# ascending order (1612 up to 2003) -- i.e. the first "part" of the new
# account-statement database, oldest period first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @(
    "1612",
    "1701", "1702", "1703", "1704", "1705", "1706", "1707", "1708", "1709", "1710", "1711", "1712",
    "1801", "1802", "1803", "1804", "1805", "1806", "1807", "1808", "1809", "1810", "1811", "1812",
    "1901", "1902", "1903", "1904", "1905", "1906", "1907", "1908", "1909", "1910", "1911", "1912",
    "2001", "2002", "2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# The data columns are "best fit" (auto-sized to content); re-fitting them
# after the text refresh nudges the stored widths slightly wider, matching
# the sizes produced when the workbook was re-saved from a newer Excel
# build. Column width is stored as characters + 5/6 (pixel) padding, so we
# back out that offset to land on the same on-disk width.
$targetWidths = @{
    2  = 18.54296875
    3  = 16.7265625
    4  = 33.90625
    5  = 13.54296875
    6  = 10.1796875
    7  = 14.36328125
    8  = 19.36328125
    9  = 18.08984375
    10 = 15
}

foreach ($col in $targetWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$col] - (5 / 6)
}
